$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.345887016753352
$ws.Range("C2").Value = 0.313515615568889
$ws.Range("D2").Value = 0.009519520851121399
$ws.Range("E2").Value = 0.046411415139886
$ws.Range("F2").Value = 4.387006073177076
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 2.666097054353429
$ws.Range("J2").Value = 0.1207758743597886
$ws.Range("L2").Value = 0.3764863800456197
$ws.Range("M2").Value = 0.5208874751191601
$ws.Range("B3").Value = 2.281277830105239
$ws.Range("C3").Value = 0.2900822958442575
$ws.Range("D3").Value = 0.00851135298239214
$ws.Range("E3").Value = 0.04606989804937722
$ws.Range("F3").Value = 4.356787725242626
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 2.653196933476522
$ws.Range("J3").Value = 0.1208226258831386
$ws.Range("L3").Value = 0.3750905753880787
$ws.Range("M3").Value = 0.5113213601050717
$ws.Range("B4").Value = 2.243105082430304
$ws.Range("C4").Value = 0.2758938340195982
$ws.Range("D4").Value = 0.007889797133564969
$ws.Range("E4").Value = 0.04585534650246803
$ws.Range("F4").Value = 4.340096582800498
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 2.646304947745094
$ws.Range("J4").Value = 0.1208522345659966
$ws.Range("L4").Value = 0.3744008067490796
$ws.Range("M4").Value = 0.5057553019092182
$ws.Range("B5").Value = 2.227925915643425
$ws.Range("C5").Value = 0.2701617096817017
$ws.Range("D5").Value = 0.007635807814867945
$ws.Range("E5").Value = 0.04576668139385109
$ws.Range("F5").Value = 4.333761774287154
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 2.6437541600136
$ws.Range("J5").Value = 0.1208645252951195
$ws.Range("L5").Value = 0.3741618199001024
$ws.Range("M5").Value = 0.5035644518145617
$ws.Range("B6").Value = 2.225428168320207
$ws.Range("C6").Value = 0.2692128896897259
$ws.Range("D6").Value = 0.007593588966432918
$ws.Range("E6").Value = 0.04575188376761119
$ws.Range("F6").Value = 4.332738043959338
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 2.643346144090216
$ws.Range("J6").Value = 0.1208665796857797
$ws.Range("L6").Value = 0.3741246802492526
$ws.Range("M6").Value = 0.5032053361788797
$ws.Range("B7").Value = 2.242898846328274
$ws.Range("C7").Value = 0.2758163275447032
$ws.Range("D7").Value = 0.007886374663616635
$ws.Range("E7").Value = 0.04585415574389984
$ws.Range("F7").Value = 4.340009260652664
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 2.646269504556926
$ws.Range("J7").Value = 0.1208523994145505
$ws.Range("L7").Value = 0.3743974131747123
$ws.Range("M7").Value = 0.5057254420653408
$ws.Range("B8").Value = 2.323298893214201
$ws.Range("C8").Value = 0.3053940953900565
$ws.Range("D8").Value = 0.009172393257621536
$ws.Range("E8").Value = 0.04629466124676185
$ws.Range("F8").Value = 4.376199182011021
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 2.661435006594331
$ws.Range("J8").Value = 0.1207918055788375
$ws.Range("L8").Value = 0.3759704138315385
$ws.Range("M8").Value = 0.5175252480542909
$ws.Range("B9").Value = 2.49286076077226
$ws.Range("C9").Value = 0.3650031772403395
$ws.Range("D9").Value = 0.0116770430952684
$ws.Range("E9").Value = 0.04712050048039718
$ws.Range("F9").Value = 4.462026099751427
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 2.69938296811074
$ws.Range("J9").Value = 0.1206802386826071
$ws.Range("L9").Value = 0.3803811184195354
$ws.Range("M9").Value = 0.5431064292208774
$ws.Range("B10").Value = 2.624731774308316
$ws.Range("C10").Value = 0.4098143013106892
$ws.Range("D10").Value = 0.01351084156434013
$ws.Range("E10").Value = 0.04770489773863495
$ws.Range("F10").Value = 4.534258423257057
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 2.732336418376107
$ws.Range("J10").Value = 0.1206028164659112
$ws.Range("L10").Value = 0.3844296206008408
$ws.Range("M10").Value = 0.5633946753979586
$ws.Range("B11").Value = 2.686317745346344
$ws.Range("C11").Value = 0.4304290492360678
$ws.Range("D11").Value = 0.01434460217728883
$ws.Range("E11").Value = 0.04796608113082357
$ws.Range("F11").Value = 4.569137325633534
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 2.748444690372992
$ws.Range("J11").Value = 0.1205686085023285
$ws.Range("L11").Value = 0.3864468418668849
$ws.Range("M11").Value = 0.5729501065158615
$ws.Range("B12").Value = 2.709868976687403
$ws.Range("C12").Value = 0.4382689693152884
$ws.Range("D12").Value = 0.0146603325096919
$ws.Range("E12").Value = 0.04806432783297154
$ws.Range("F12").Value = 4.582637475784765
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 2.754706332899161
$ws.Range("J12").Value = 0.1205558026801987
$ws.Range("L12").Value = 0.3872359401601955
$ws.Range("M12").Value = 0.5766154731498006
$ws.Range("B13").Value = 2.704786565057248
$ws.Range("C13").Value = 0.4365790015310722
$ws.Range("D13").Value = 0.01459233311491914
$ws.Range("E13").Value = 0.04804319770104826
$ws.Range("F13").Value = 4.579716949586611
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 2.753350563159572
$ws.Range("J13").Value = 0.1205585540209619
$ws.Range("L13").Value = 0.3870648723605399
$ws.Range("M13").Value = 0.5758239830789122
$ws.Range("B14").Value = 2.688250708778924
$ws.Range("C14").Value = 0.4310733686547223
$ws.Range("D14").Value = 0.01437057705015832
$ws.Range("E14").Value = 0.04797417705637219
$ws.Range("F14").Value = 4.570242123348379
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 2.748956590733457
$ws.Range("J14").Value = 0.1205675519788452
$ws.Range("L14").Value = 0.3865112561200164
$ws.Range("M14").Value = 0.5732507176793931
$ws.Range("B15").Value = 2.678151976234631
$ws.Range("C15").Value = 0.4277053952940264
$ws.Range("D15").Value = 0.01423474748310127
$ws.Range("E15").Value = 0.04793181463203755
$ws.Range("F15").Value = 4.564476633811466
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 2.746286257819804
$ws.Range("J15").Value = 0.1205730828359335
$ws.Range("L15").Value = 0.3861754340227179
$ws.Range("M15").Value = 0.5716806311213389
$ws.Range("B16").Value = 2.620739149363658
$ws.Range("C16").Value = 0.4084717483320901
$ws.Range("D16").Value = 0.01345635014524049
$ws.Range("E16").Value = 0.04768773605483112
$ws.Range("F16").Value = 4.532019805758296
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 2.731306285836766
$ws.Range("J16").Value = 0.1206050726096439
$ws.Range("L16").Value = 0.3843013208000627
$ws.Range("M16").Value = 0.5627767715149687
$ws.Range("B17").Value = 2.585927508899147
$ws.Range("C17").Value = 0.3967317730724744
$ws.Range("D17").Value = 0.0129787625812412
$ws.Range("E17").Value = 0.0475368166620127
$ws.Range("F17").Value = 4.512627241172794
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 2.722403564886775
$ws.Range("J17").Value = 0.1206249582104229
$ws.Range("L17").Value = 0.3831965561835915
$ws.Range("M17").Value = 0.5573981062948192
$ws.Range("B18").Value = 2.566055089040617
$ws.Range("C18").Value = 0.3900008627708758
$ws.Range("D18").Value = 0.01270402106389668
$ws.Range("E18").Value = 0.04744957226347424
$ws.Range("F18").Value = 4.501663180193333
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 2.717388082092938
$ws.Range("J18").Value = 0.1206364909652589
$ws.Range("L18").Value = 0.3825776483508747
$ws.Range("M18").Value = 0.5543351425186884
$ws.Range("B19").Value = 2.559352437766449
$ws.Range("C19").Value = 0.387725589271497
$ws.Range("D19").Value = 0.0126109884093637
$ws.Range("E19").Value = 0.04741995693378254
$ws.Range("F19").Value = 4.497983521173126
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 2.715707948592936
$ws.Range("J19").Value = 0.1206404120125826
$ws.Range("L19").Value = 0.3823709354275877
$ws.Range("M19").Value = 0.5533033481178933
$ws.Range("B20").Value = 2.589617708472645
$ws.Range("C20").Value = 0.3979792730550002
$ws.Range("D20").Value = 0.01302960684535037
$ws.Range("E20").Value = 0.04755292769063502
$ws.Range("F20").Value = 4.514671932615784
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 2.723340386696407
$ws.Range("J20").Value = 0.1206228314952149
$ws.Range("L20").Value = 0.3833124503025118
$ws.Range("M20").Value = 0.5579674966760706
$ws.Range("B21").Value = 2.693101447043432
$ws.Range("C21").Value = 0.4326895916016156
$ws.Range("D21").Value = 0.01443571155624568
$ws.Range("E21").Value = 0.04799446784389971
$ws.Range("F21").Value = 4.573017163112809
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 2.750242807345714
$ws.Range("J21").Value = 0.1205649050244619
$ws.Range("L21").Value = 0.3866731823171108
$ws.Range("M21").Value = 0.574005274148945
$ws.Range("B22").Value = 2.762074601580537
$ws.Range("C22").Value = 0.4555706570098437
$ws.Range("D22").Value = 0.01535473282939392
$ws.Range("E22").Value = 0.04827921174058325
$ws.Range("F22").Value = 4.612853542821341
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 2.768768637299601
$ws.Range("J22").Value = 0.1205279105167931
$ws.Range("L22").Value = 0.3890166102817574
$ws.Range("M22").Value = 0.5847604620902374
$ws.Range("B23").Value = 2.725139565077768
$ws.Range("C23").Value = 0.4433405131272821
$ws.Range("D23").Value = 0.01486420698550717
$ws.Range("E23").Value = 0.04812758464713074
$ws.Range("F23").Value = 4.591435569027396
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 2.758794358560564
$ws.Range("J23").Value = 0.1205475753107654
$ws.Range("L23").Value = 0.3877524351300679
$ws.Range("M23").Value = 0.5789951740951622
$ws.Range("B24").Value = 2.587948929285346
$ws.Range("C24").Value = 0.3974152205466908
$ws.Range("D24").Value = 0.01300662068814518
$ws.Range("E24").Value = 0.04754564538702155
$ws.Range("F24").Value = 4.513746951569146
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 2.722916529407613
$ws.Range("J24").Value = 0.1206237926712261
$ws.Range("L24").Value = 0.3832600039978189
$ws.Range("M24").Value = 0.5577099840156379
$ws.Range("B25").Value = 2.445711987294374
$ws.Range("C25").Value = 0.348701539680178
$ws.Range("D25").Value = 0.0110008872604368
$ws.Range("E25").Value = 0.0469010725628225
$ws.Range("F25").Value = 4.437204897777747
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 2.688231168493758
$ws.Range("J25").Value = 0.1207096310574554
$ws.Range("L25").Value = 0.3790459967362381
$ws.Range("M25").Value = 0.5359241212781356
